# Update gh-pages output values (column F: "想去人数") across sheets
$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 47
$ws.Range("F3").Value = 148
$ws.Range("F4").Value = 146
$ws.Range("F5").Value = 2052
$ws.Range("F6").Value = 4190
$ws.Range("F7").Value = 543
$ws.Range("F8").Value = 1050
$ws.Range("F9").Value = 661
$ws.Range("F10").Value = 374
$ws.Range("F12").Value = 2203
$ws.Range("F14").Value = 659379
$ws.Range("F15").Value = 1627
$ws.Range("F16").Value = 517
$ws.Range("F17").Value = 1458
$ws.Range("F19").Value = 542
$ws.Range("F20").Value = 1273
$ws.Range("F21").Value = 2237
$ws.Range("F22").Value = 1137
$ws.Range("F23").Value = 2695
$ws.Range("F24").Value = 1549
$ws.Range("F25").Value = 809
$ws.Range("F26").Value = 1535
$ws.Range("F27").Value = 26
$ws.Range("F28").Value = 528
$ws.Range("F29").Value = 1085
$ws.Range("F30").Value = 279
$ws.Range("F31").Value = 1082
$ws.Range("F33").Value = 79
$ws.Range("F34").Value = 2019
$ws.Range("F35").Value = 1369
$ws.Range("F37").Value = 1239
$ws.Range("F38").Value = 2503
$ws.Range("F39").Value = 1144
$ws.Range("F40").Value = 29
$ws.Range("F42").Value = 2579
$ws.Range("F43").Value = 209
$ws.Range("F44").Value = 983
$ws.Range("F45").Value = 3125
$ws.Range("F48").Value = 878
$ws.Range("F49").Value = 154
$ws.Range("F50").Value = 670

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 51
$ws.Range("F9").Value = 105
$ws.Range("F10").Value = 479
$ws.Range("F11").Value = 144616
$ws.Range("F12").Value = 144616
$ws.Range("F13").Value = 10
$ws.Range("F18").Value = 230
$ws.Range("F19").Value = 333
$ws.Range("F23").Value = 125
$ws.Range("F24").Value = 84
$ws.Range("F27").Value = 558
$ws.Range("F32").Value = 347
$ws.Range("F33").Value = 272
$ws.Range("F38").Value = 194
$ws.Range("F41").Value = 187

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 3130
$ws.Range("F5").Value = 243
$ws.Range("F7").Value = 824
$ws.Range("F8").Value = 1180
$ws.Range("F10").Value = 1594
$ws.Range("F12").Value = 91
$ws.Range("F13").Value = 1893

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 824
$ws.Range("F3").Value = 1180
$ws.Range("F5").Value = 148
$ws.Range("F6").Value = 1594
$ws.Range("F8").Value = 146
$ws.Range("F9").Value = 2052
$ws.Range("F10").Value = 91
$ws.Range("F11").Value = 1893
$ws.Range("F12").Value = 4190
$ws.Range("F13").Value = 543
$ws.Range("F14").Value = 661
$ws.Range("F15").Value = 374
$ws.Range("F16").Value = 2203
$ws.Range("F18").Value = 659387
$ws.Range("F19").Value = 105
$ws.Range("F20").Value = 479
$ws.Range("F21").Value = 1627
$ws.Range("F22").Value = 144616
$ws.Range("F23").Value = 1458
$ws.Range("F25").Value = 542
$ws.Range("F26").Value = 1273
$ws.Range("F27").Value = 2237
$ws.Range("F28").Value = 1137
$ws.Range("F29").Value = 2695
$ws.Range("F30").Value = 1549
$ws.Range("F31").Value = 809
$ws.Range("F33").Value = 1535
$ws.Range("F35").Value = 528
$ws.Range("F36").Value = 125
$ws.Range("F37").Value = 1085
$ws.Range("F38").Value = 1082
$ws.Range("F39").Value = 79
$ws.Range("F40").Value = 2019
$ws.Range("F41").Value = 1369
$ws.Range("F42").Value = 1239
$ws.Range("F43").Value = 2503
$ws.Range("F44").Value = 1144
$ws.Range("F45").Value = 347
$ws.Range("F46").Value = 347
$ws.Range("F47").Value = 272
$ws.Range("F48").Value = 2579
$ws.Range("F49").Value = 209
$ws.Range("F50").Value = 983
$ws.Range("F51").Value = 3125
$ws.Range("F52").Value = 154
$ws.Range("F53").Value = 670
